$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("I2").Value = 0.6275843652054141
$ws.Range("J2").Value = 0.6275843652054141
$ws.Range("M2").Value = 0.073876
$ws.Range("S2").Value = 0.6275843652054141
$ws.Range("T2").Value = 0.6275843652054141

# Row 3 updates
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.008832333333333333
$ws.Range("H3").Value = 0.026497
$ws.Range("I3").Value = 0.372415634794586
$ws.Range("J3").Value = 0.3724156347945861
$ws.Range("M3").Value = 0.073876
$ws.Range("Q3").Value = 0.0006524974573333332
$ws.Range("R3").Value = 0.005872477116
$ws.Range("S3").Value = 0.372415634794586
$ws.Range("T3").Value = 0.3724156347945861
